$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 274 (shifts IAD..YHZ down by one row)
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row with the NQZ / Astana, Kazakhstan entry
$ws.Cells.Item(274, 1).Value = "NQZ"
$ws.Cells.Item(274, 2).Value = "ASTANA, Kazakhstan"
$ws.Cells.Item(274, 3).Value = "Asia Pacific"
$ws.Cells.Item(274, 4).Value = "ASTANA"
$ws.Cells.Item(274, 5).Value = "Kazakhstan"
$ws.Cells.Item(274, 6).Value = "KZ"
$ws.Cells.Item(274, 7).Value = 51.167801
$ws.Cells.Item(274, 8).Value = 71.418893

# Row-insert already carries over the bold font + center/top alignment used
# by the rest of column A; restore the thin box border to fully match the
# style used by every other colo-code cell.
$ws.Cells.Item(274, 1).Borders.LineStyle = 1
